{"js": "// Update the date line and the 5x5 grid of \"two-digit \u00f7 one-digit\" answers.\n// Cells are addressed by (row, col) position rather than by searching for\n// their old text, because some old values reappear elsewhere as NEW values\n// (e.g. \"67\u00f78=8, 3\" is the old text of one cell and the new text of\n// another), which would make a blind global find/replace ambiguous.\n\nconst body = context.document.body;\n\n// 1) Title paragraph: \"2025-03-18 Tuesday\" -> \"2025-03-19 Wednesday\"\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(\"2025-03-19 Wednesday\", \"Replace\");\n\n// 2) Table of answers (5 data rows among the 20 table rows; 5 columns).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// rowIndex -> [col0, col1, col2, col3, col4] new values\nconst newValues = {\n  0: [\"95\u00f78=11, 7\", \"57\u00f75=11, 2\", \"76\u00f73=25, 1\", \"80\u00f76=13, 2\", \"76\u00f77=10, 6\"],\n  4: [\"99\u00f73=33, 0\", \"87\u00f73=29, 0\", \"64\u00f72=32, 0\", \"23\u00f75=4, 3\", \"12\u00f72=6, 0\"],\n  8: [\"62\u00f73=20, 2\", \"31\u00f74=7, 3\", \"26\u00f79=2, 8\", \"47\u00f78=5, 7\", \"93\u00f74=23, 1\"],\n  12: [\"73\u00f78=9, 1\", \"89\u00f74=22, 1\", \"52\u00f78=6, 4\", \"63\u00f78=7, 7\", \"67\u00f78=8, 3\"],\n  16: [\"92\u00f74=23, 0\", \"41\u00f77=5, 6\", \"33\u00f77=4, 5\", \"60\u00f77=8, 4\", \"34\u00f78=4, 2\"],\n};\n\nfor (const rowIndexStr of Object.keys(newValues)) {\n  const rowIndex = Number(rowIndexStr);\n  const rowValues = newValues[rowIndex];\n  for (let col = 0; col < rowValues.length; col++) {\n    table.getCell(rowIndex, col).value = rowValues[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 5x5 grid of \"two-digit \u00f7 one-digit\" answers.\n# Cells are addressed by (row, col) position rather than by searching for\n# their old text, because some old values reappear elsewhere as NEW values\n# (e.g. \"67\u00f78=8, 3\" is the old text of one cell and the new text of\n# another), which would make a blind global find/replace ambiguous.\n\n$d = $word.ActiveDocument\n\n# 1) Title paragraph: \"2025-03-18 Tuesday\" -> \"2025-03-19 Wednesday\"\n$d.Paragraphs.Item(1).Range.Text = \"2025-03-19 Wednesday\"\n\n# 2) Table of answers (5 data rows among the 20 table rows; 5 columns).\n#    Word's Cell(row, col) is 1-based.\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text = \"95\u00f78=11, 7\"\n$t.Cell(1, 2).Range.Text = \"57\u00f75=11, 2\"\n$t.Cell(1, 3).Range.Text = \"76\u00f73=25, 1\"\n$t.Cell(1, 4).Range.Text = \"80\u00f76=13, 2\"\n$t.Cell(1, 5).Range.Text = \"76\u00f77=10, 6\"\n\n$t.Cell(5, 1).Range.Text = \"99\u00f73=33, 0\"\n$t.Cell(5, 2).Range.Text = \"87\u00f73=29, 0\"\n$t.Cell(5, 3).Range.Text = \"64\u00f72=32, 0\"\n$t.Cell(5, 4).Range.Text = \"23\u00f75=4, 3\"\n$t.Cell(5, 5).Range.Text = \"12\u00f72=6, 0\"\n\n$t.Cell(9, 1).Range.Text = \"62\u00f73=20, 2\"\n$t.Cell(9, 2).Range.Text = \"31\u00f74=7, 3\"\n$t.Cell(9, 3).Range.Text = \"26\u00f79=2, 8\"\n$t.Cell(9, 4).Range.Text = \"47\u00f78=5, 7\"\n$t.Cell(9, 5).Range.Text = \"93\u00f74=23, 1\"\n\n$t.Cell(13, 1).Range.Text = \"73\u00f78=9, 1\"\n$t.Cell(13, 2).Range.Text = \"89\u00f74=22, 1\"\n$t.Cell(13, 3).Range.Text = \"52\u00f78=6, 4\"\n$t.Cell(13, 4).Range.Text = \"63\u00f78=7, 7\"\n$t.Cell(13, 5).Range.Text = \"67\u00f78=8, 3\"\n\n$t.Cell(17, 1).Range.Text = \"92\u00f74=23, 0\"\n$t.Cell(17, 2).Range.Text = \"41\u00f77=5, 6\"\n$t.Cell(17, 3).Range.Text = \"33\u00f77=4, 5\"\n$t.Cell(17, 4).Range.Text = \"60\u00f77=8, 4\"\n$t.Cell(17, 5).Range.Text = \"34\u00f78=4, 2\"\n"}
